# Criação das classes AdicionarLinhas, EmpurrarLinhasParaBaixo e atualização
# da classe EscritorExcelModelo: insere novas linhas de item no orçamento,
# empurrando a linha do TOTAL para baixo, e atualiza o total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EmpurrarLinhasParaBaixo: abre espaço para 4 novas linhas de item logo
# após a última linha existente (linha 9), empurrando a linha em branco
# e a linha do TOTAL (antiga linha 11) para baixo (nova linha 15).
$ws.Rows("10:13").Insert()

# AdicionarLinhas: replica a formatação/estilo da última linha de item
# (linha 9 - "Limpeza do sistema") para as quatro novas linhas.
$ws.Range("A9:D9").Copy()
$ws.Range("A10").PasteSpecial()
$ws.Range("A9:D9").Copy()
$ws.Range("A11").PasteSpecial()
$ws.Range("A9:D9").Copy()
$ws.Range("A12").PasteSpecial()
$ws.Range("A9:D9").Copy()
$ws.Range("A13").PasteSpecial()

# Valores das novas linhas (Quantidade / Valor Unitário / Subtotal).
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 300
$ws.Range("D10").Value = 300

$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 100
$ws.Range("D11").Value = 100

$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 100
$ws.Range("D12").Value = 100

$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 100
$ws.Range("D13").Value = 100

# Atualiza o TOTAL (agora na linha 15) com a soma recalculada.
$ws.Range("D15").Value = 1300
